$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.311.21"

$ws.Range("D3").Value = "1.873.57"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'0.7118"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'241.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.3110"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").Value = "'0.07779"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.99%  "

$ws.Range("D10").Value = "'25.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("D11").Value = "'0.08387"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "1.867.58"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").Value = "'5.241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").Value = "'91.15"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "29.319.76"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "'6.061"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("D18").Value = "'0.000008199"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.14%  "

$ws.Range("D19").Value = "'240.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").Value = "2.119.08"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'7.758"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'0.1586"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "'162.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").Value = "'9.020"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").Value = "'1.511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "'1.290"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").Value = "'4.318"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "

$ws.Range("D33").Value = "'0.05297"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.99%  "

$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("E35").Value = "  +1.34%  "

$ws.Range("D36").Value = "'0.7432"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.64%  "

$ws.Range("D37").Value = "'2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").Value = "'0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "

$ws.Range("D39").Value = "1.228.95"
$ws.Range("E39").Value = "  +5.51%  "

$ws.Range("D40").Value = "'2.731"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").Value = "'6.528"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.73%  "

$ws.Range("D42").Value = "'0.8842"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").Value = "'109.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.37%  "

$ws.Range("D44").Value = "'72.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").Value = "2.016.24"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").Value = "'1.796"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5193"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000123"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "

$ws.Range("D50").Value = "'9.386"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").Value = "'0.4310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
